$d = $word.ActiveDocument

# --- Part 1: adjust table column widths in the first table (3305 -> 3304 dxa, 3831 -> 3832 dxa) ---
# Cell.Width is expressed in points (1 pt = 20 dxa); Word propagates the change
# to every row of the fixed-layout column and to the table grid definition.
$t = $d.Tables(1)
$t.Cell(1, 2).Width = 165.2
$t.Cell(1, 3).Width = 191.6

# --- Part 2: append new paragraphs documenting the "Implementazione del tempo" notes ---
$n = $d.Paragraphs.Count
$r = $d.Paragraphs($n).Range

$r.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$r = $d.Paragraphs($n).Range
$r.Text = 'Implementazione del tempo:'

$r.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$r = $d.Paragraphs($n).Range
$r.Text = 'voglio implementare la gestione del tempo in questa maniera:'

$r.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$r = $d.Paragraphs($n).Range
$r.Text = 'ho una interfaccia azione che viene implementata da due classi una azioni sincrone ed un’altra azioni asincrone.'

$r.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$r = $d.Paragraphs($n).Range
$r.Text = 'Poi le azioni dell’attaccante saranno azioniSincrone ed azioniAsincrone dove le azioni sincrone dell’attaccante e del difensore andranno ad interagire con una variabile Tempo che avrà un delta iniziale di 120 ovvero un dominio intero di [+60;-60].'

$r.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$r = $d.Paragraphs($n).Range
$r.Text = 'Quando il difensore deve eseguire una mossa andrà a verificare che il tempo sia <=0, se ciò è vero lui può eseguire una mossa e fare tempo + X, dove X sarà positivo.'

$r.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$r = $d.Paragraphs($n).Range
$r.Text = 'L’attaccante quando andrà ad eseguire una mossa verificherà che il tempo sia >=0, se ciò è vero lui potrà eseguire la sua mossa e fare tempo +Y con Y negativo.'

$r.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$r = $d.Paragraphs($n).Range
$r.Text = 'Ciò cosa rappresenta nella gestione del tempo, ovvero se il difensore trova un tempo negativo vuol dire che il tempo sta andando avanti e l’attaccante sta facendo, e lui può agire.'

$r.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$r = $d.Paragraphs($n).Range
$r.Text = 'Viceversa se l’attaccante trova un tempo positivo vuol dire che il tempo va avanti ed il difensore sta facendo un’azione.'

$r.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$r = $d.Paragraphs($n).Range
$r.Text = 'Ovvero per esempio: se l’attaccante trova il tempo negativo vuol dire che lui stesso sta facendo un’azione che ancora non termina e che se non termina non può fare altro (sia ipotizzando un fattore di risorse o più semplicemente perché posso fare un exploit solo dopo aver concluso un portscan e scoperto i servizi attaccabili).'

$r.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$r = $d.Paragraphs($n).Range
$r.Text = 'Invece le azioni asincrone, sono azioni che non vincolano le altre ma permettono la ‘parallelizzazione’ temporale ma non istantanea. Ovvero il difensore può sempre fare un’azione per turno, in un turno non può scegliere più azioni asincrone da fare insieme, ma al suo turno successivo può eseguire un’altra azione senza attendere che termini l’altra precedentemente attivata.'

$r.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$r = $d.Paragraphs($n).Range
$r.Text = 'L’intenzione è quella di dare vita ad un nuovo agente che attende la terminazione del tempo per quell’azione, mentre l’agente principale continua ad operare nelle operazioni successive.'

$r.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$r = $d.Paragraphs($n).Range
$r.Text = 'Solo allo scadere del tempo l’agente asincrono modifica lo stato e si distrugge.'

$r.InsertParagraphAfter()
$n = $d.Paragraphs.Count
$r = $d.Paragraphs($n).Range
$r.Text = 'Capire se è proprio necessario definire un vero agente nell’ambiente o basta lanciare un sottoprocesso e tutto funziona (capire bene quelle asincrone come implementarle)'
